# Updated spreadsheets with new tests.
#
# Adds 16 new "Test_Cases" rows to Sheet1 (rows 130-145), updates the JIRA
# Issue Number for two existing rows (111/112), and refreshes the sheet
# selection to match the new extent of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the JIRA Issue Number on the two existing InputValueCalc
#     error-diagnostic rows (DFDL-148 -> DFDL-280) ---------------------
$ws.Range("H111").Value = "DFDL-280"
$ws.Range("H112").Value = "DFDL-280"

# --- New test case rows (columns: B,C,D,E,F,G,H,I) --------------------
# Column meaning: B=Test Case Name, C=Status, D=Key Requirement ID,
# E=Priority, F=Spin Version, G=Input Data Type, H=JIRA Issue Number,
# I=Negative Test?
$newRows = @(
    @("byte_01",                     "Passed", "DFDL-5-015R",  "High", 2, "text", "DFDL-184", "Yes"),
    @("inputValueCalcAbsolutePath",  "Passed", "DFDL-17-007R", "High", 4, "text", "DFDL-283", $null),
    @("int_error_03",                "Passed", "DFDL-5-013R",  "High", 4, "text", "DFDL-279", "Yes"),
    @("short_02",                    "Passed", "DFDL-5-014R",  "High", 4, "text", "DFDL-279", "Yes"),
    @("unsignedInt_02",              "Passed", "DFDL-5-018R",  "High", 4, "text", "DFDL-279", "Yes"),
    @("byte_02",                     "Passed", "DFDL-5-015R",  "High", 4, "text", "DFDL-279", "Yes"),
    @("unsignedByte_02",             "Passed", "DFDL-5-020R",  "High", 4, "text", "DFDL-279", "Yes"),
    @("lke1_rel",                    "Passed", "DFDL-23-011R", "High", 4, "byte", "DFDL-237", $null),
    @("lke1_abs",                    "Passed", "DFDL-23-011R", "High", 4, "byte", "DFDL-237", $null),
    @("ocke1",                       "Passed", "DFDL-23-011R", "High", 4, "byte", "DFDL-239", $null),
    @("ocke2",                       "Passed", "DFDL-23-011R", "High", 4, "byte", "DFDL-239", $null),
    @("InputValueCalc_01",           "Passed", "DFDL-17-007R", "High", 4, "text", "DFDL-236", $null),
    @("InputValueCalc_02",           "Passed", "DFDL-17-007R", "High", 4, "text", "DFDL-236", "Yes"),
    @("InputValueCalc_03",           "Passed", "DFDL-17-007R", "High", 4, "text", "DFDL-236", "Yes"),
    @("InputValueCalc_05",           "Passed", "DFDL-17-007R", "High", 4, "text", "DFDL-236", $null),
    @("InputValueCalc_06",           "Passed", "DFDL-17-007R", "High", 4, "text", "DFDL-236", $null)
)

$startRow = 130
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $ws.Range("G$r").Value = $vals[5]
    $ws.Range("H$r").Value = $vals[6]
    if ($vals[7] -ne $null) {
        $ws.Range("I$r").Value = $vals[7]
    }
}

# --- Refresh selection to reflect the new bottom of the table ---------
$ws.Activate() | Out-Null
$ws.Range("H145").Select() | Out-Null
